# Convert HOUR_APPR_PROCESS_START values (numeric hours) into a time-like
# text string of the form "HH:00:00", e.g. 14 -> "14:00:00".
#
# The HOUR_APPR_PROCESS_START column is column V (22nd column), with data
# rows starting at row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerRow = 1
$hourCol = 22  # column V

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt $headerRow) {
    $lastRow = $headerRow
}

for ($r = $headerRow + 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $hourCol)
    $hour = $cell.Value2
    if ($null -ne $hour -and "$hour" -ne "") {
        $cell.Value = "$($hour):00:00"
    }
}
